# "massive MDY site template update"
# Insert three new columns (Month/Day/Year) before the existing "Date Sampled"
# column, mirroring the style of the neighbouring Transect (D) column, and
# populate them from the already-present Date Sampled serial values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Make room: insert 3 blank columns at E:G (old E "Date Sampled" etc. shift to H:Q)
$ws.Columns("E:G").Insert()

# 2. Header row (row 1): copy the bold/centered header format already used by
#    B1/C1/D1 onto the three new header cells, then give them their labels.
$ws.Range("B1").Copy()
$ws.Range("E1:G1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("E1").Value = "Month"
$ws.Range("F1").Value = "Day"
$ws.Range("G1").Value = "Year"

# 3. Data rows: copy the Transect column's (D) data-cell format onto E2:G5,
#    then fill in Month/Day/Year derived from each row's Date Sampled value.
$ws.Range("D2").Copy()
$ws.Range("E2:G5").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("E2").Value = 8
$ws.Range("F2").Value = 13
$ws.Range("G2").Value = 2015

$ws.Range("E3").Value = 8
$ws.Range("F3").Value = 13
$ws.Range("G3").Value = 2015

$ws.Range("E4").Value = 8
$ws.Range("F4").Value = 11
$ws.Range("G4").Value = 2015

$ws.Range("E5").Value = 8
$ws.Range("F5").Value = 12
$ws.Range("G5").Value = 2015

# 4. Widen the new columns to match column D's width.
$ws.Range("E1:G1").ColumnWidth = 17.33203125

# 5. Leave the selection where the edit ended.
$ws.Range("E6").Select()

Write-Output "done"
